$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row23 = @('6422273','Lemon Juice','1','16.55','16.55')
$row24 = @('8255796','Tuna White Chunk (Pouch)','6','72.00','432.00')
$row25 = @('1132582','Sour Cream','1','28.94','28.94')
$row26 = @('7722184','Parmesan (Grated)','1','59.95','59.95')
$row27 = @('1027629','Cheddar - (Sliced)','12','34.91','418.92')
$row28 = @('1035842','Feta - Pail','1','92.87','92.87')
$row29 = @('6364494','Yogurt - Greek (Bulk)','2','27.11','54.22')
$row30 = @('4254710','Veggie Burger','1','52.72','52.72')
$row31 = @('1365278','Vegan Chicken Tenders','1','87.80','87.80')
$row32 = @('3275539','Sauerkraut','3','19.35','58.05')
$row33 = @('1028188','Tortellini - Cheese','2','32.22','64.44')
$row34 = @('7529232','Wrap - Wheat (10")','2','31.32','62.64')
$row35 = @('2825368','Sausage - Chicken Patty','3','50.01','150.03')
$row36 = @('4157160','Spanakopita','2','75.51','151.02')
$row37 = @('9546982','Arugula - Fresh','4','20.35','81.40')
$row38 = @('6264014','Tomato - Grape','6','18.06','108.36')

$newRows = @($row23, $row24, $row25, $row26, $row27, $row28, $row29, $row30, $row31, $row32, $row33, $row34, $row35, $row36, $row37, $row38)

$startRow = 23
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 0; $c -lt 5; $c++) {
        $cell = $ws.Cells.Item($r, $c + 1)
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$c]
    }
}
